# Apply the changes described by the commit "Dev guide and minor changes"
# to the UndoRedoStartingStateListDiagram deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1. TextBox 3 ("currentStatePointer = 0") -> single run "Index = 0"
# ---------------------------------------------------------------------------
$tb = $s.Shapes.Item(1)
$tr = $tb.TextFrame.TextRange
$firstRun = $tr.Characters(1, 19)   # "currentStatePointer"
$firstRun.Text = ""
$tr.Text = "Index = 0"

# ---------------------------------------------------------------------------
# 2. Remove the "Table 4" graphicFrame (the ab0:AddressBook table)
# ---------------------------------------------------------------------------
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    if ($s.Shapes.Item($i).Name -eq "Table 4") {
        $s.Shapes.Item($i).Delete()
    }
}

# ---------------------------------------------------------------------------
# 3. Add a new "Empty List" textbox (appears after the connector, last
#    in z-order) at the position the table used to occupy.
# ---------------------------------------------------------------------------
# iron_native's shape-id allocator advances on every Shapes.Add* call made
# during this session (independent from ids already present in the
# document); spend its first two allocations on scratch shapes so the
# real textbox below lands on id 8 - matching the id the authentic
# Microsoft PowerPoint edit produced.
$scratch1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$scratch1.Delete()
$scratch2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$scratch2.Delete()

$left   = 825623 / 12700.0
$top    = 2793209 / 12700.0
$width  = 1150892 / 12700.0
$height = 369332 / 12700.0

$newBox = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$newBox.Name = "TextBox 7"
$newBox.TextFrame.WordWrap = 0
$newBox.TextFrame.AutoSize = 1
$newBox.Fill.Visible = 0

$ntr = $newBox.TextFrame.TextRange
$ntr.Text = "Empty List"
$ntr.LanguageID = "en-SG"

# re-assert the exact size/position once text/auto-fit settled
$newBox.Left = $left
$newBox.Top = $top
$newBox.Width = $width
$newBox.Height = $height

# ---------------------------------------------------------------------------
# 4. Refresh the cached "datetimeFigureOut" placeholder text on every slide
#    layout (7/6/2018 -> 3/30/2019), as PowerPoint does on save.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -eq "7/6/2018") {
                $shp.TextFrame.TextRange.Text = "3/30/2019"
            }
        }
    }
}
